# Commit: "add rural communities interactions back in, fix de_dg files (not run yet)"
#
# The stakeholder-interaction matrix on Sheet1 had an extra row for
# "Legislature" (row 12) that needs to be removed; every row below it
# shifts up by one. Once no cell references the "Legislature" shared
# string anymore, Excel drops it from the shared-strings table and
# reindexes the remaining strings, which is exactly what the diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Select the entire "Legislature" row, then delete it (rows below shift up).
$ws.Rows(12).Select()
$ws.Rows(12).Delete()

# After the delete, Excel leaves the same on-screen row selected (now holding
# what used to be row 13, "Friant Water Authority"), so reproduce that
# resulting selection state.
$ws.Range("A12:XFD12").Select()

# Best-effort: also nudge the scrolled viewport to roughly match the
# post-edit view (top-left visible cell around row 10).
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
